$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.139.93'
$ws.Range('E2').Value = '  +1.78%  '
$ws.Range('D3').Value = '3.901.61'
$ws.Range('E3').Value = '  +2.86%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '467.58'
$ws.Range('E5').Value = '  +9.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.59'
$ws.Range('E6').Value = '  +3.86%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.628'
$ws.Range('E7').Value = '  +0.94%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.165'
$ws.Range('E10').Value = '  +7.86%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000340'
$ws.Range('E11').Value = '  +8.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.08'
$ws.Range('E12').Value = '  +0.30%  '
$ws.Range('E13').Value = '  -0.97%  '
$ws.Range('D14').Value = '4.529.05'
$ws.Range('E14').Value = '  +3.27%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.39'
$ws.Range('E15').Value = '  +2.07%  '
$ws.Range('D16').Value = '3.929.08'
$ws.Range('E16').Value = '  +1.62%  '
$ws.Range('E17').Value = '  -0.36%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.96'
$ws.Range('E18').Value = '  -0.26%  '
$ws.Range('E19').Value = '  +2.63%  '
$ws.Range('D20').Value = '67.351.42'
$ws.Range('E20').Value = '  +1.73%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '433.16'
$ws.Range('E21').Value = '  +6.66%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.77'
$ws.Range('E22').Value = '  -2.77%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.36'
$ws.Range('E23').Value = '  +3.58%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.81'
$ws.Range('E24').Value = '  +4.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '38.87'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.54'
$ws.Range('E26').Value = '  +7.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '5.73'
$ws.Range('E27').Value = '  +6.09%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.10'
$ws.Range('E28').Value = '  +1.77%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.61'
$ws.Range('E29').Value = '  -3.38%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '737.47'
$ws.Range('E30').Value = '  +4.79%  '
$ws.Range('E31').Value = '  -1.50%  '
$ws.Range('E32').Value = '  +0.37%  '
$ws.Range('E33').Value = '  +0.51%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '43.25'
$ws.Range('E34').Value = '  +6.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.158'
$ws.Range('E35').Value = '  +4.76%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '58.15'
$ws.Range('E36').Value = '  +2.89%  '
$ws.Range('E37').Value = '  -0.12%  '
$ws.Range('D38').Value = '0.0₃0793'
$ws.Range('E38').Value = '  +16.77%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.40'
$ws.Range('E39').Value = '  -6.59%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.26'
$ws.Range('E40').Value = '  +13.86%  '
$ws.Range('E41').Value = '  +0.17%  '
$ws.Range('E42').Value = '  -0.79%  '
$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.337'
$ws.Range('E43').Value = '  +5.04%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.80'
$ws.Range('E45').Value = '  +5.49%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.17'
$ws.Range('E46').Value = '  +5.03%  '
$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.50'
$ws.Range('E47').Value = '  -5.67%  '
$ws.Range('B48').Value = 'LidoDAOToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.41'
$ws.Range('E48').Value = '  +1.07%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.91'
$ws.Range('E49').Value = '  +3.21%  '
$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.14'
$ws.Range('E50').Value = '  -0.22%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '143.51'
$ws.Range('E51').Value = '  +0.82%  '
